# Apply updated simulated-game probability values to the team matrix sheet.
# (added more games, sped up simulate game logic, and drafted optimization logic)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B2" = 0.1739130434782609
    "C2" = 0.5217391304347826
    "P2" = 0.1304347826086956
    "S2" = 0.1739130434782609
    "P3" = 0.5
    "S3" = 0.5
    "P4" = 0.5
    "S4" = 0.5
    "B6" = 0.03225806451612903
    "D6" = 0.03225806451612903
    "J6" = 0.3225806451612903
    "O6" = 0.06451612903225806
    "Q6" = 0.0967741935483871
    "R6" = 0.06451612903225806
    "S6" = 0.3870967741935484
    "B7" = 0.03703703703703703
    "J7" = 0.03703703703703703
    "O7" = 0.03703703703703703
    "Q7" = 0.1851851851851852
    "R7" = 0.07407407407407407
    "S7" = 0.6296296296296297
    "B8" = 0.0847457627118644
    "F8" = 0.05084745762711865
    "J8" = 0.05084745762711865
    "Q8" = 0.06779661016949153
    "R8" = 0.1355932203389831
    "S8" = 0.6101694915254238
    "F9" = 0.1333333333333333
    "J9" = 0.06666666666666667
    "O9" = 0.06666666666666667
    "Q9" = 0.1333333333333333
    "R9" = 0.06666666666666667
    "S9" = 0.5333333333333333
    "B10" = 0.0763888888888889
    "D10" = 0.006944444444444444
    "F10" = 0.1041666666666667
    "J10" = 0.09027777777777778
    "O10" = 0.006944444444444444
    "Q10" = 0.1736111111111111
    "R10" = 0.09027777777777778
    "S10" = 0.4513888888888889
    "G11" = 0.1458333333333333
    "J11" = 0.1041666666666667
    "K11" = 0.2291666666666667
    "L11" = 0.5
    "S11" = 0.02083333333333333
    "G12" = 0.64
    "J12" = 0.32
    "S12" = 0.04
    "F13" = 0.1111111111111111
    "G13" = 0.5555555555555556
    "F15" = 0.1379310344827586
    "H15" = 0.103448275862069
    "I15" = 0.103448275862069
    "J15" = 0.3448275862068966
    "K15" = 0.03448275862068965
    "S15" = 0.2758620689655172
    "H16" = 0.2
    "I16" = 0.1
    "J16" = 0.3
    "K16" = 0.1
    "O16" = 0.2
    "S16" = 0.1
    "H17" = 0.125
    "I17" = 0.025
    "J17" = 0.425
    "K17" = 0.225
    "O17" = 0.15
    "S17" = 0.05
    "H18" = 0.4615384615384616
    "I18" = 0.03846153846153846
    "J18" = 0.1153846153846154
    "K18" = 0.2307692307692308
    "M18" = 0.03846153846153846
    "O18" = 0.03846153846153846
    "S18" = 0.07692307692307693
    "F19" = 0.02762430939226519
    "H19" = 0.2099447513812155
    "I19" = 0.04972375690607735
    "J19" = 0.3812154696132597
    "K19" = 0.1104972375690608
    "M19" = 0.03867403314917127
    "O19" = 0.05524861878453038
    "S19" = 0.1270718232044199
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
